# Working on Optimize project test
#
# Update the "optimize" sheet's key/value table with new optimization
# parameters, and leave the UI focused on that sheet (cell B9 selected),
# matching the Excel session that produced the commit.

$wb = $excel.ActiveWorkbook

# Leftover selection state on the "data" sheet from before the user
# navigated away to work on "optimize".
$wsData = $wb.Worksheets.Item("data")
$wsData.Range("D1:E1048576").Select() | Out-Null

$ws = $wb.Worksheets.Item("optimize")

# ending_regularization: 1.000000 -> 10
$ws.Range("B2").Value = 10

# number_of_particles: 128 -> 32
$ws.Range("B8").Value = 32

# procrustes: false -> TRUE
$ws.Range("B10").Value = $true

# procrustes_interval: 0 -> 1
$ws.Range("B11").Value = 1

# procrustes_scaling: false -> TRUE
$ws.Range("B12").Value = $true

# starting_regularization: 10.000000 -> 1000
$ws.Range("B14").Value = 1000

# Switch focus to the optimize sheet, with B9 selected.
$ws.Activate() | Out-Null
$ws.Range("B9").Select() | Out-Null
